$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New s_val data (regenerated to filter save games)
$data = @{
    2 = @(0.7287194209349384, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 1, 3.034748368925986)
    3 = @(1.505614041169197, 9.226618575922256, 0.1529057820181812, 6.48142807727062, 1, 17.36656647638025)
    4 = @(1.505614041169197, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 1, 6.741336633845642)
    5 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 6.48142807727062, 0, 12.0302756157461)
    6 = @(0.7287194209349384, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 1, 3.034748368925986)
    7 = @(0.7287194209349384, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 0, 3.034748368925986)
    8 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 0, 6.048734245549538)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
    $ws.Cells.Item($r, 6).Value = $vals[4]
    $ws.Cells.Item($r, 7).Value = $vals[5]
}
